$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 40: D40 formula gains "+225000"
$ws.Range("D40").Formula = "=45000+225000"

# Row 41: C41 formula gains "+15000000+59269000"
$ws.Range("C41").Formula = "=37292000+2308000+15000000+59269000"

# Row 42: D42 formula gains more terms
$ws.Range("D42").Formula = "=39600000+4800000+1266000+3915000+1935000+1014000+15000000+2100000"

# Row 43 (new): SALES - cash/retail
$ws.Range("B43").Value = "SALES - cash/retail"
$ws.Range("C43").Formula = "=55531525+14163475-59269000"

# Row 44 (new): SELISIH - lebih
$ws.Range("B44").Value = "SELISIH - lebih"
$ws.Range("C44").Value = 23000

# Row 45 (new): SETOR KE BANK
$ws.Range("B45").Value = "SETOR KE BANK"
$ws.Range("D45").Value = 55000000

# Row 46 (new): Wages Expense on 16-Jan-2021 (serial date 44212)
$ws.Range("A46").Value = 44212
$ws.Range("B46").Value = "Wages Expense"

# Apply the same number formats as the analogous existing cells (s="1" / s="2")
$ws.Range("C43").NumberFormat = $ws.Range("C41").NumberFormat
$ws.Range("C44").NumberFormat = $ws.Range("C41").NumberFormat
$ws.Range("D45").NumberFormat = $ws.Range("D42").NumberFormat
$ws.Range("A46").NumberFormat = $ws.Range("A40").NumberFormat

# Update the selection to reflect the new view position (keep the existing
# freeze at row 2/col 0 intact; just move the active cell to C45)
$ws.Activate()
$win = $excel.ActiveWindow
if (-not $win.FreezePanes) {
  $ws.Range("A3").Select()
  $win.FreezePanes = $true
}
$ws.Range("C45").Select()
